$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header/numeric values and shared-string text values, preserving styles.
$ws.Range("A3").Value = 122
$ws.Range("B3").Value = "TES-098"

$ws.Range("A4").Value = 321
$ws.Range("B4").Value = "TES-123"

$ws.Range("A5").Value = 1234
$ws.Range("B5").Value = "123-CON"

$ws.Range("A6").Value = 3030301
$ws.Range("B6").Value = "123-09"

# C2 text stays "No. Buku" (same text, but was rewritten by shared-string reshuffle)
$ws.Range("C2").Value = "No. Buku"

# Move active selection to D2
$ws.Range("D2").Select()
